# "Added hours for week 3"
# Row 6 on the "Time Sheet" worksheet is Week 3. Hyrum's (column B/C) hours
# for week 3 were entered (B6=9), which ripples the running Total column
# (C) down through the rest of the sheet. At the same time, Patrick's
# Total column (G) is switched from a plain "=F{n}" restatement of that
# week's hours to a proper running total "=F{n}+G{n-1}", matching the
# pattern already used by the other members' Total columns (C, E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New hours entry for Hyrum, week 3 (row 6).
$ws.Range("B6").Value = 9

# Patrick's weekly Total column (G) becomes a running total of column F,
# same pattern as columns C and E. G4 (the first week) is left untouched.
$ws.Range("G5").Formula  = "=F5+G4"
$ws.Range("G6").Formula  = "=F6+G5"
$ws.Range("G7").Formula  = "=F7+G6"
$ws.Range("G8").Formula  = "=F8+G7"
$ws.Range("G9").Formula  = "=F9+G8"
$ws.Range("G10").Formula = "=F10+G9"
$ws.Range("G11").Formula = "=F11+G10"
$ws.Range("G12").Formula = "=F12+G11"
$ws.Range("G13").Formula = "=F13+G12"

# Cosmetic row-height tweak that came along with the edit.
$ws.Rows.Item(12).RowHeight = 15

# Update the view: scroll so row 4 is the top visible row, and move the
# selection to G14.
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("G14").Select() | Out-Null
